$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 314-315, shifting existing rows 314:429 down to 316:431
$ws.Rows("314:315").Insert()

# --- New row 314 ---
$ws.Cells.Item(314, 1).Value2 = 10
$ws.Cells.Item(314, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(314, 3).Value2 = "La Araucanía"
$ws.Cells.Item(314, 4).Value2 = 44809
$ws.Cells.Item(314, 5).Value2 = 9
$ws.Cells.Item(314, 6).Value2 = 100112037
$ws.Cells.Item(314, 7).Value2 = "Cebollín"
$ws.Cells.Item(314, 8).Value2 = "Sin especificar"
$ws.Cells.Item(314, 9).Value2 = "Primera"
$ws.Cells.Item(314, 10).Value2 = 70
$ws.Cells.Item(314, 11).Value2 = 8000
$ws.Cells.Item(314, 12).Value2 = 9000
$ws.Cells.Item(314, 13).Value2 = 8429
$ws.Cells.Item(314, 14).Value2 = "`$/docena de paquetes"
$ws.Cells.Item(314, 15).Value2 = "Provincia de Cautín"
$ws.Cells.Item(314, 16).Value2 = 702
$ws.Cells.Item(314, 17).Value2 = 12
$ws.Cells.Item(314, 18).Value2 = "Hortaliza"

# --- New row 315 ---
$ws.Cells.Item(315, 1).Value2 = 10
$ws.Cells.Item(315, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(315, 3).Value2 = "La Araucanía"
$ws.Cells.Item(315, 4).Value2 = 44809
$ws.Cells.Item(315, 5).Value2 = 9
$ws.Cells.Item(315, 6).Value2 = 100112037
$ws.Cells.Item(315, 7).Value2 = "Cebollín"
$ws.Cells.Item(315, 8).Value2 = "Sin especificar"
$ws.Cells.Item(315, 9).Value2 = "Primera"
$ws.Cells.Item(315, 10).Value2 = 100
$ws.Cells.Item(315, 11).Value2 = 7000
$ws.Cells.Item(315, 12).Value2 = 7000
$ws.Cells.Item(315, 13).Value2 = 7000
$ws.Cells.Item(315, 14).Value2 = "`$/docena de paquetes"
$ws.Cells.Item(315, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(315, 16).Value2 = 583
$ws.Cells.Item(315, 17).Value2 = 12
$ws.Cells.Item(315, 18).Value2 = "Hortaliza"
